# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts from the site refresh at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rId1 / sheet1) ---
$wsExhibit = $wb.Worksheets.Item("展览")

$exhibitUpdates = @{
    2  = 0
    4  = 41
    5  = 0
    9  = 200
    10 = 1290
    11 = 0
    12 = 0
    14 = 0
    15 = 17
    17 = 46
    18 = 9
    19 = 5043
    20 = 0
    22 = 0
    24 = 0
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# --- Sheet "全部类型" (rId4 / sheet4) ---
$wsAll = $wb.Worksheets.Item("全部类型")

$allTypesUpdates = @{
    2  = 0
    3  = 94
    4  = 41
    5  = 0
    6  = 154
    7  = 0
    8  = 66
    10 = 0
    13 = 400
    14 = 140
    17 = 46
    18 = 9
    20 = 5043
    21 = 0
    23 = 0
    24 = 474
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}
